$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value corrections (addressed using the ORIGINAL row
# numbers, before the row deletions below take place) ---

# D3 (RM 8): was missing -> -14.2
$ws.Range("D3").Value = -14.2

# F4 (RM 9): was 17.97 -> missing
$ws.Range("F4").Value = "'"
$ws.Range("F4").Style = "Normal"

# D5 (RM 14): was -14.4 -> missing
$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"

# F9 (RM 42): was missing -> 17.26
$ws.Range("F9").Value = 17.26

# F10 (RM 52 a): was missing -> 16.43
$ws.Range("F10").Value = 16.43

# F17 (RM 116): was 17.78 -> missing
$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"

# F18 (RM 120): was 18.35 -> missing
$ws.Range("F18").Value = "'"
$ws.Range("F18").Style = "Normal"

# D21 (RM 135): was missing -> -14.3
$ws.Range("D21").Value = -14.3

# D23 (RM 140): was -13.9 -> missing
$ws.Range("D23").Value = "'"
$ws.Range("D23").Style = "Normal"

# D34 (SC 193): was missing -> -14.7
$ws.Range("D34").Value = -14.7

# --- Whole-row deletions (bottom row first so the earlier row index used
# for the second delete stays valid) ---

# Row 28 ("SC 92") removed entirely; rows below shift up.
$ws.Rows.Item(28).Delete()

# Row 26 ("RM 232") removed entirely; rows below shift up.
$ws.Rows.Item(26).Delete()
